$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rewrite A1 title without the trailing newline character.
#    (Content is identical to before, just strips the trailing "\n";
#    Excel will naturally move this string to the end of the shared
#    strings table since the old entry becomes unused.)
$title = "8.3.1.2 Экономикадагы иш менен камсыз болгон бардык калктын чакан жана орто ишканаларда иштегендердин үлүшү"
$ws.Range("A1").Value = $title

# 2) Add the new "2023" column (N) of data, matching the formatting of
#    the existing last data column (M).
$ws.Range("M3:M6").Copy()
$ws.Range("N3:N6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(4,14).Value = 2023
$ws.Cells.Item(5,14).Value = 2.5449890821474286
$ws.Cells.Item(6,14).Value = 1.4569686017619159

# 3) Minor row-height tweaks that accompanied the new column.
$ws.Rows.Item(1).RowHeight = 45
$ws.Rows.Item(5).RowHeight = 17.25
$ws.Rows.Item(6).RowHeight = 17.25

# 4) Reset the stored selection back to the top-left cell.
$ws.Range("A1").Select() | Out-Null
